$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp string
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 08:16"

# Alemania (row 18): active cases / recovered updated
$ws.Range("D18").Value = 181000
$ws.Range("E18").Value = 6653

# Kazajistan moves above Portugal (row 39/40 swap with refreshed Kazajistan figures)
$ws.Range("A39").Value = "Kazajistan"
$ws.Range("B39").Value = 44075
$ws.Range("C39").Value = 1501
$ws.Range("D39").Value = 26251
$ws.Range("E39").Value = 17636
$ws.Range("H39").Value = 188

$ws.Range("A40").Value = "Portugal"
$ws.Range("B40").Value = 42782
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 28097
$ws.Range("E40").Value = 13098
$ws.Range("H40").Value = 1587

# Rumania (row 50): active cases / recovered updated
$ws.Range("D50").Value = 20433
$ws.Range("E50").Value = 5626

# Israel moves above Nigeria (row 51/52 swap with refreshed Israel figures)
$ws.Range("A51").Value = "Israel"
$ws.Range("B51").Value = 27542
$ws.Range("C51").Value = 495
$ws.Range("D51").Value = 17599
$ws.Range("E51").Value = 9618
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 325

$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 27110
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 10801
$ws.Range("E52").Value = 15693
$ws.Range("H52").Value = 616

# Uzbekistan (row 72): totals / new cases / recovered updated
$ws.Range("B72").Value = 9199
$ws.Range("C72").Value = 121
$ws.Range("E72").Value = 3138

# Pacific / Caribbean islands reordered (values unchanged, only names shuffled)
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
